# adding HoH gender to strata variables
$wb = $excel.ActiveWorkbook

$wsGeneral = $wb.Worksheets.Item("general")
$wsVariables = $wb.Worksheets.Item("variables")
$wsStrata = $wb.Worksheets.Item("strata_variables")

# ---------------------------------------------------------------
# Sheet "general": filename lost its "_2" duplicate suffix (1 of 2)
# ---------------------------------------------------------------
$wsGeneral.Range("B3").Value = "BFA2402_MSNA_2024_DATA_CLEANED_MV.xlsx"

# ---------------------------------------------------------------
# Sheet "strata_variables": add the new "HoH gender" variable name
# used by each country/source in the strata configuration
# ---------------------------------------------------------------
$wsStrata.Range("D3").Value = "c_chef_menage_genre_final"   # BFA
$wsStrata.Range("H4").Value = "hoh_gender_final"             # CAR
$wsStrata.Range("K5").Value = "hoh_gender_final"             # DRC
$wsStrata.Range("J6").Value = "hoh_gender_final"             # ETH
$wsStrata.Range("K9").Value = "HHhGenderFinal"                # KEN1
$wsStrata.Range("I10").Value = "HHhGenderFinal"               # KEN2
$wsStrata.Range("H11").Value = "hoh_gender_final"             # MLI

# ---------------------------------------------------------------
# Sheet "general": second filename lost its "_2" duplicate suffix
# ---------------------------------------------------------------
$wsGeneral.Range("B13").Value = "REACH_MSNA-2024-NIGER_Base-de-donnees_Septembre2024.xlsx"

# ---------------------------------------------------------------
# Sheet "strata_variables": remaining HoH gender additions
# ---------------------------------------------------------------
$wsStrata.Range("I13").Value = "c_gender_hoh"                 # NER
$wsStrata.Range("H14").Value = "final_hoh_gender"             # SOM
$wsStrata.Range("J16").Value = "Q4_4_hoh_gender"              # SYR

# ---------------------------------------------------------------
# Sheet "variables": the helper/reference columns (B:J) get hidden,
# and the stale style on AB6:AE6 is cleared back to Normal
# ---------------------------------------------------------------
$wsVariables.Range("B1:J1").EntireColumn.Hidden = $true
$wsVariables.Range("AB6:AE6").Style = "Normal"

# ---------------------------------------------------------------
# Selection / active-sheet bookkeeping: restore each sheet's last
# selection, finishing on strata_variables so it becomes the
# active (visible) tab, matching the saved workbook state
# ---------------------------------------------------------------
$wsGeneral.Range("B14").Select()
$wsVariables.Range("S30").Select()
$wsStrata.Range("O17").Select()
$wsStrata.Activate()
